# Replace text values like "0,1" (comma decimal separator) in column C
# with true numeric values using a period decimal separator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 0.1
    3  = 0.1
    4  = 0.1
    5  = 0.1
    6  = 0.1
    7  = 0.1
    8  = 0.1
    9  = 0.1
    10 = 0.1
    11 = 0.1
    12 = 0.25
    13 = 0.2
    14 = 0.2
    15 = 0.1
    16 = 0.1
    17 = 0.1
    18 = 0.05
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}

# Match the new active-cell selection recorded in the saved sheet view.
$ws.Range("J10").Select()
